# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# per-language report sheets ("zh-cn" and "de-de"), row 2
# (the row for eaa5a729-29d1-4618-bfde-a68e29ca8271).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 07:16:35"
$wsZhCn.Range("H2").Value = "2016-03-23 07:17:11"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 07:16:42"
$wsDeDe.Range("H2").Value = "2016-03-23 07:17:24"
